# paises.xlsx ("Pais" sheet) refresh: COVID-19 country stats pulled at 17:29
# (previous pull was at 16:12), plus four countries (Chile, Republica de
# Macedonia, Mali, Benin) moving a few rows earlier in the ranking as their
# case counts overtook their neighbours.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: refreshed timestamp banner
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 17 de Julio de 2020 a las 17:29"

# Peru/Sudafrica/Mexico block: Chile moves up to 9, pushing Sudafrica -> 10, Mexico -> 11
$ws.Cells.Item(9, 1).Value = "Chile"
$ws.Cells.Item(9, 2).Value = 326539
$ws.Cells.Item(9, 3).Value = 2841
$ws.Cells.Item(9, 4).Value = 296814
$ws.Cells.Item(9, 5).Value = 22435
$ws.Cells.Item(9, 8).Value = 7290
$ws.Cells.Item(10, 1).Value = "Sudafrica"
$ws.Cells.Item(10, 2).Value = 324221
$ws.Cells.Item(10, 3).Value = 0
$ws.Cells.Item(10, 4).Value = 165591
$ws.Cells.Item(10, 5).Value = 153961
$ws.Cells.Item(10, 7).Value = 0
$ws.Cells.Item(10, 8).Value = 4669
$ws.Cells.Item(11, 1).Value = "Mexico"
$ws.Cells.Item(11, 2).Value = 324041
$ws.Cells.Item(11, 3).Value = 6406
$ws.Cells.Item(11, 4).Value = 203464
$ws.Cells.Item(11, 5).Value = 83003
$ws.Cells.Item(11, 7).Value = 668
$ws.Cells.Item(11, 8).Value = 37574

# Noruega block: Republica de Macedonia moves up to 80, pushing Malasia -> 81
$ws.Cells.Item(80, 1).Value = "Republica de Macedonia"
$ws.Cells.Item(80, 2).Value = 8786
$ws.Cells.Item(80, 3).Value = 163
$ws.Cells.Item(80, 4).Value = 4676
$ws.Cells.Item(80, 5).Value = 3704
$ws.Cells.Item(80, 7).Value = 5
$ws.Cells.Item(80, 8).Value = 406
$ws.Cells.Item(81, 1).Value = "Malasia"
$ws.Cells.Item(81, 2).Value = 8755
$ws.Cells.Item(81, 3).Value = 18
$ws.Cells.Item(81, 4).Value = 8541
$ws.Cells.Item(81, 5).Value = 92
$ws.Cells.Item(81, 8).Value = 122

# Libano block: Mali moves up to 114, pushing Cuba -> 115
$ws.Cells.Item(114, 1).Value = "Mali"
$ws.Cells.Item(114, 2).Value = 2467
$ws.Cells.Item(114, 3).Value = 27
$ws.Cells.Item(114, 4).Value = 1791
$ws.Cells.Item(114, 5).Value = 555
$ws.Cells.Item(114, 8).Value = 121
$ws.Cells.Item(115, 1).Value = "Cuba"
$ws.Cells.Item(115, 2).Value = 2444
$ws.Cells.Item(115, 3).Value = 4
$ws.Cells.Item(115, 4).Value = 2300
$ws.Cells.Item(115, 5).Value = 57
$ws.Cells.Item(115, 8).Value = 87

# Libia block: Benin moves up to 129, pushing Suazilandia/Yemen/Nueva Zelanda/Ruanda down one row each
$ws.Cells.Item(129, 1).Value = "Benin"
$ws.Cells.Item(129, 2).Value = 1602
$ws.Cells.Item(129, 3).Value = 139
$ws.Cells.Item(129, 4).Value = 782
$ws.Cells.Item(129, 5).Value = 789
$ws.Cells.Item(129, 7).Value = 3
$ws.Cells.Item(129, 8).Value = 31
$ws.Cells.Item(130, 1).Value = "Suazilandia"
$ws.Cells.Item(130, 4).Value = 736
$ws.Cells.Item(130, 5).Value = 795
$ws.Cells.Item(130, 8).Value = 21
$ws.Cells.Item(131, 1).Value = "Yemen"
$ws.Cells.Item(131, 2).Value = 1552
$ws.Cells.Item(131, 3).Value = 0
$ws.Cells.Item(131, 4).Value = 695
$ws.Cells.Item(131, 5).Value = 419
$ws.Cells.Item(131, 8).Value = 438
$ws.Cells.Item(132, 1).Value = "Nueva Zelanda"
$ws.Cells.Item(132, 2).Value = 1549
$ws.Cells.Item(132, 3).Value = 1
$ws.Cells.Item(132, 4).Value = 1506
$ws.Cells.Item(132, 5).Value = 21
$ws.Cells.Item(132, 8).Value = 22
$ws.Cells.Item(133, 1).Value = "Ruanda"
$ws.Cells.Item(133, 2).Value = 1473
$ws.Cells.Item(133, 4).Value = 770
$ws.Cells.Item(133, 5).Value = 699
$ws.Cells.Item(133, 8).Value = 4

# Remaining rows: statistics refreshed in place (no reordering)
$ws.Cells.Item(4, 2).Value = 3717343
$ws.Cells.Item(4, 3).Value = 22318
$ws.Cells.Item(4, 4).Value = 1681017
$ws.Cells.Item(4, 5).Value = 1894896
$ws.Cells.Item(4, 7).Value = 312
$ws.Cells.Item(4, 8).Value = 141430
$ws.Cells.Item(6, 2).Value = 1020644
$ws.Cells.Item(6, 3).Value = 15007
$ws.Cells.Item(6, 5).Value = 350695
$ws.Cells.Item(13, 2).Value = 293239
$ws.Cells.Item(13, 3).Value = 687
$ws.Cells.Item(13, 7).Value = 114
$ws.Cells.Item(13, 8).Value = 45233
$ws.Cells.Item(26, 2).Value = 88171
$ws.Cells.Item(26, 3).Value = 2023
$ws.Cells.Item(26, 4).Value = 56495
$ws.Cells.Item(26, 5).Value = 28060
$ws.Cells.Item(26, 7).Value = 94
$ws.Cells.Item(26, 8).Value = 3616
$ws.Cells.Item(43, 2).Value = 50113
$ws.Cells.Item(43, 3).Value = 1370
$ws.Cells.Item(43, 4).Value = 24423
$ws.Cells.Item(43, 5).Value = 24748
$ws.Cells.Item(43, 7).Value = 1
$ws.Cells.Item(43, 8).Value = 942
$ws.Cells.Item(44, 2).Value = 48077
$ws.Cells.Item(44, 3).Value = 312
$ws.Cells.Item(44, 4).Value = 32790
$ws.Cells.Item(44, 5).Value = 13605
$ws.Cells.Item(44, 7).Value = 3
$ws.Cells.Item(44, 8).Value = 1682
$ws.Cells.Item(56, 2).Value = 26636
$ws.Cells.Item(56, 3).Value = 471
$ws.Cells.Item(56, 4).Value = 17805
$ws.Cells.Item(56, 5).Value = 8490
$ws.Cells.Item(56, 7).Value = 7
$ws.Cells.Item(56, 8).Value = 341
$ws.Cells.Item(61, 2).Value = 20494
$ws.Cells.Item(61, 3).Value = 230
$ws.Cells.Item(61, 4).Value = 13913
$ws.Cells.Item(61, 5).Value = 5906
$ws.Cells.Item(61, 7).Value = 9
$ws.Cells.Item(61, 8).Value = 675
$ws.Cells.Item(63, 2).Value = 19439
$ws.Cells.Item(63, 3).Value = 169
$ws.Cells.Item(63, 4).Value = 17335
$ws.Cells.Item(63, 5).Value = 1393
$ws.Cells.Item(90, 2).Value = 6786
$ws.Cells.Item(90, 3).Value = 45
$ws.Cells.Item(90, 4).Value = 5483
$ws.Cells.Item(90, 5).Value = 1247
$ws.Cells.Item(110, 2).Value = 2782
$ws.Cells.Item(110, 3).Value = 4
$ws.Cells.Item(110, 4).Value = 2591
$ws.Cells.Item(110, 5).Value = 154
$ws.Cells.Item(112, 2).Value = 2689
$ws.Cells.Item(112, 3).Value = 2
$ws.Cells.Item(112, 5).Value = 666
$ws.Cells.Item(134, 2).Value = 1402
$ws.Cells.Item(134, 3).Value = 19
$ws.Cells.Item(134, 4).Value = 397
$ws.Cells.Item(134, 5).Value = 996
$ws.Cells.Item(136, 2).Value = 1336
$ws.Cells.Item(136, 3).Value = 9
$ws.Cells.Item(136, 4).Value = 1095
$ws.Cells.Item(136, 5).Value = 191
$ws.Cells.Item(156, 2).Value = 614
$ws.Cells.Item(156, 3).Value = 2
$ws.Cells.Item(156, 5).Value = 139
$ws.Cells.Item(187, 2).Value = 85
$ws.Cells.Item(187, 3).Value = 1
$ws.Cells.Item(187, 5).Value = 3
